# Updates crypto price/volume data (and two coin-row swaps) to match
# the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.224.61"
$ws.Range("E2").Value = "  +4.95%  "
$ws.Range("D3").Value = "4.075.13"
$ws.Range("E3").Value = "  +5.21%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.96"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.86"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.717"
$ws.Range("E7").Value = "  +18.12%  "
$ws.Range("D8").Value = "4.065.71"
$ws.Range("E8").Value = "  +5.13%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.777"
$ws.Range("E10").Value = "  +8.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.182"
$ws.Range("E11").Value = "  +7.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000340"
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.97"
$ws.Range("E13").Value = "  +16.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.09"
$ws.Range("E14").Value = "  +8.27%  "
$ws.Range("D15").Value = "4.725.17"
$ws.Range("E15").Value = "  +5.42%  "
$ws.Range("D16").Value = "4.090.43"
$ws.Range("E16").Value = "  +5.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.49"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "21.31"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D21").Value = "72.279.03"
$ws.Range("E21").Value = "  +5.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.56"
$ws.Range("E22").Value = "  +6.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.55"
$ws.Range("E23").Value = "  +18.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.63"
$ws.Range("E24").Value = "  +6.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.02"
$ws.Range("E25").Value = "  +6.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.08"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.22"
$ws.Range("E28").Value = "  +5.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.00"
$ws.Range("E29").Value = "  +5.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.83"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("E31").Value = "  +16.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.72"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +5.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "686.19"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.79"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.68"
$ws.Range("E36").Value = "  +12.77%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0906"
$ws.Range("E37").Value = "  +6.65%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "42.62"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +5.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  +8.30%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  +12.30%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.98"
$ws.Range("E47").Value = "  +17.15%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000293"
$ws.Range("E49").Value = "  +8.56%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.39"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +4.58%  "

Write-Output "Applied cryptos update"
